# Updated location, added dinner
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the day headers (column A) from March 2024 dates to November 2024 dates
$ws.Range("A2").Value = "25/11/2024"
$ws.Range("A8").Value = "26/11/2024"
$ws.Range("A15").Value = "27/11/2024"
$ws.Range("A19").Value = "28/11/2024"
$ws.Range("A25").Value = "29/11/2024"

# Rename "Dinner" to "Course Dinner"
$ws.Range("E24").Value = "Course Dinner"

# Update the selection to match the saved workbook state
$ws.Range("E32").Select()
